$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 23, mirroring the formatting of row 22 (the prior last row),
# with the next day's expense data.
$ws.Range("A22:M22").Copy()
$ws.Range("A23:M23").PasteSpecial(-4122)

$ws.Range("A23").Value = 43811
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 16.5
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 3

# Update the selection to match the new active cell.
$ws.Range("G23").Select()
